# Update "想去人数" (interested-count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for sheet "展览"
$updatesSheet1 = @{
    6  = 76
    7  = 4392
    14 = 1666
    15 = 668
    18 = 334
    24 = 28
    26 = 564
    32 = 1077
    35 = 1169
    37 = 281
    39 = 553
    45 = 114
    48 = 74
}

foreach ($row in $updatesSheet1.Keys) {
    $sheet1.Range("F$row").Value = $updatesSheet1[$row]
}

# Row -> new F value, for sheet "全部类型"
$updatesSheet4 = @{
    5  = 76
    6  = 4392
    9  = 1666
    12 = 668
    15 = 334
    20 = 28
    22 = 564
    30 = 1077
    35 = 281
    40 = 553
    46 = 114
    48 = 74
}

foreach ($row in $updatesSheet4.Keys) {
    $sheet4.Range("F$row").Value = $updatesSheet4[$row]
}
